$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.508.06'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +5.46%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.791.06'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +23.04%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '617.47'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +8.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.19'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.33%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.776.88'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +22.66%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.548'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +6.77%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.171'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +13.55%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.41'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.34%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.505'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +8.53%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.74'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +13.81%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000262'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +9.63%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.450.80'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +23.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.808.68'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +23.56%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '70.730.30'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +5.89%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.61'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +9.25%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '525.06'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +8.69%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.00'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.42'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +22.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.748'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +9.49%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.12'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +5.85%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +11.94%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.58'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +7.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.92'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +6.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000124'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +34.47%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.51'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +9.81%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +13.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.97'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.97%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '32.35'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +16.19%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +4.33%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.30%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.17'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +11.73%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.05'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +11.25%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +10.60%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.19'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +10.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.133'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +8.52%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '51.67'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +5.64%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.157.71'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +13.26%  '
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '430.57'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +16.94%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.89'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +8.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '44.43'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.83%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +4.62%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0369'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +8.13%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.75'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +8.28%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '140.65'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +4.21%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +11.15%  '
